$d = $word.ActiveDocument

# 1) Mark the inline picture's run as "no proofing" (w:noProof) so Word
#    doesn't re-check spelling/grammar on the field-like lastRenderedPageBreak
#    run that carries the drawing.
$shape = $d.InlineShapes(1)
$shape.Range.NoProofing = $true

# 2) Remove the trailing scratch notes (segment truth table + PC7...PC0
#    byte encodings) that followed the table, now that PC0...PC7 are
#    documented as GPIO outputs. Keep the blank paragraph right after the
#    table.
$count = $d.Paragraphs.Count
$firstToDelete = $null
for ($i = 1; $i -le $count; $i++) {
    $t = $d.Paragraphs($i).Range.Text
    if ($t.StartsWith("0: a, b, c, d, e, f --> 1  g --> 0")) {
        $firstToDelete = $i
        break
    }
}

if ($firstToDelete -ne $null) {
    $startPara = $d.Paragraphs($firstToDelete)
    $lastPara = $d.Paragraphs($d.Paragraphs.Count)
    $r = $d.Range($startPara.Range.Start, $lastPara.Range.End)
    $r.Delete()
}
